# Add an "Implemented?" column before the existing "Finalized?" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the existing "Finalized?" header (and its formatting) from D to E,
# then insert the new "Implemented?" header in D.
$ws.Range("E1").Value = $ws.Range("D1").Value2
$ws.Range("D1").Value = "Implemented?"

# Copy the "Bad" style used in column D (rows 2-7) over to the new column E.
$ws.Range("E2:E7").Style = $ws.Range("D2:D7").Style

# Match the new column's width to the target sheet (closest value this
# interop's character-width rounding can produce to 18.109375).
$ws.Range("E1").ColumnWidth = 17.3

# Update the selection to mirror the recorded end-user state.
$ws.Range("D22").Select()
